$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value, preventing Excel from
# auto-converting numeric-looking strings into numbers, while leaving
# the cell style unchanged afterwards.
function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "30.787.16"
$ws.Range("E2").Value = "  +0.97%  "

# Row 3
$ws.Range("D3").Value = "1.963.82"
$ws.Range("E3").Value = "  +3.88%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.9889"
$ws.Range("E4").Value = "  -1.19%  "

# Row 5
Set-TextValue $ws.Range("D5") "252.64"
$ws.Range("E5").Value = "  +3.57%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.7001"
$ws.Range("E6").Value = "  +48.83%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.9939"
$ws.Range("E7").Value = "  -0.66%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3269"
$ws.Range("E8").Value = "  +12.68%  "

# Row 9
Set-TextValue $ws.Range("D9") "25.84"
$ws.Range("E9").Value = "  +16.36%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.06890"
$ws.Range("E10").Value = "  +5.99%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.8397"
$ws.Range("E11").Value = "  +15.57%  "

# Row 12
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D12") "101.49"
$ws.Range("E12").Value = "  +5.91%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D13") "0.07957"
$ws.Range("E13").Value = "  +2.63%  "

# Row 14
$ws.Range("D14").Value = "1.942.79"
$ws.Range("E14").Value = "  +2.79%  "

# Row 15
Set-TextValue $ws.Range("D15") "5.391"
$ws.Range("E15").Value = "  +3.78%  "

# Row 16
Set-TextValue $ws.Range("D16") "278.35"
$ws.Range("E16").Value = "  -1.18%  "

# Row 17
$ws.Range("D17").Value = "30.814.57"
$ws.Range("E17").Value = "  +1.07%  "

# Row 18
Set-TextValue $ws.Range("D18") "13.86"
$ws.Range("E18").Value = "  +6.07%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.000007701"
$ws.Range("E19").Value = "  +3.07%  "

# Row 20
Set-TextValue $ws.Range("D20") "5.635"
$ws.Range("E20").Value = "  +6.54%  "

# Row 21
$ws.Range("D21").Value = "2.196.60"
$ws.Range("E21").Value = "  +2.82%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.9969"
$ws.Range("E22").Value = "  -0.39%  "

# Row 23
Set-TextValue $ws.Range("D23") "0.9948"
$ws.Range("E23").Value = "  -0.54%  "

# Row 24
Set-TextValue $ws.Range("D24") "6.659"
$ws.Range("E24").Value = "  +6.25%  "

# Row 25
Set-TextValue $ws.Range("D25") "9.538"
$ws.Range("E25").Value = "  +4.99%  "

# Row 26
Set-TextValue $ws.Range("D26") "165.20"
$ws.Range("E26").Value = "  +0.70%  "

# Row 27
Set-TextValue $ws.Range("D27") "19.55"
$ws.Range("E27").Value = "  +3.22%  "

# Row 28
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D28") "0.1284"
$ws.Range("E28").Value = "  +31.93%  "

# Row 29
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D29") "2.152"
$ws.Range("E29").Value = "  +13.32%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D30") "1.356"
$ws.Range("E30").Value = "  +1.87%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "1.555"
$ws.Range("E31").Value = "  +5.79%  "

# Row 32
Set-TextValue $ws.Range("D32") "4.512"
$ws.Range("E32").Value = "  +5.31%  "

# Row 33
Set-TextValue $ws.Range("D33") "4.370"
$ws.Range("E33").Value = "  +5.29%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.05055"
$ws.Range("E34").Value = "  +4.06%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.208"
$ws.Range("E35").Value = "  +7.13%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.7314"
$ws.Range("E36").Value = "  +5.08%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.699"
$ws.Range("E37").Value = "  -0.67%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.01981"
$ws.Range("E38").Value = "  +4.53%  "

# Row 39
Set-TextValue $ws.Range("D39") "2.939"
$ws.Range("E39").Value = "  +3.91%  "

# Row 40
Set-TextValue $ws.Range("D40") "6.555"
$ws.Range("E40").Value = "  +5.38%  "

# Row 41
Set-TextValue $ws.Range("D41") "77.68"
$ws.Range("E41").Value = "  +3.00%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.4662"
$ws.Range("E42").Value = "  +9.64%  "

# Row 43
Set-TextValue $ws.Range("D43") "2.042"
$ws.Range("E43").Value = "  +2.42%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.8467"
$ws.Range("E44").Value = "  +2.70%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.9942"
$ws.Range("E45").Value = "  -0.60%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "10.03"
$ws.Range("E46").Value = "  +4.16%  "

# Row 47
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D47") "102.34"
$ws.Range("E47").Value = "  +0.87%  "

# Row 48
Set-TextValue $ws.Range("D48") "7.380"
$ws.Range("E48").Value = "  +5.93%  "

# Row 49
Set-TextValue $ws.Range("D49") "36.11"
$ws.Range("E49").Value = "  +2.85%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.4226"
$ws.Range("E50").Value = "  +7.31%  "

# Row 51
Set-TextValue $ws.Range("D51") "932.54"
$ws.Range("E51").Value = "  +2.12%  "
